$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ratios")

# Update column T (20) width: was 23.7109375, now 22.7109375 (in stored OOXML units).
# Excel's ColumnWidth uses character-width units that get pixel-snapped on save;
# 21.83 is the input value that yields the closest achievable stored width to the target.
$ws.Columns.Item(20).ColumnWidth = 21.83

# Updated computed ratio values (columns P:U, rows 2-12) after removing a
# redundant slope constant from the calculation.
$ws.Range("P2").Value = 1.311545646351963
$ws.Range("Q2").Value = 0.5314624646829403
$ws.Range("R2").Value = 0.17118120129988
$ws.Range("S2").Value = 1.612128312904708
$ws.Range("T2").Value = 0.2257985283109158
$ws.Range("U2").Value = 1.726512590317885
$ws.Range("P3").Value = 0.03522926649515275
$ws.Range("Q3").Value = 5.693201124441734
$ws.Range("R3").Value = 0.3296613420644465
$ws.Range("S3").Value = 1.210058348418021
$ws.Range("T3").Value = 0.01157263513266091
$ws.Range("U3").Value = 2.453369603464564
$ws.Range("P4").Value = 1.311366584465857
$ws.Range("Q4").Value = 0.695716825952851
$ws.Range("R4").Value = 0.172954589408702
$ws.Range("S4").Value = 1.285989101452288
$ws.Range("T4").Value = 0.2271469555391469
$ws.Range("U4").Value = 1.004207115162461
$ws.Range("P5").Value = 0.005919331541588404
$ws.Range("Q5").Value = 2.812994683903398
$ws.Range("R5").Value = 0.2260885021375802
$ws.Range("S5").Value = 1.381858217885171
$ws.Range("T5").Value = 0.001274687244363013
$ws.Range("U5").Value = 2.977480699240536
$ws.Range("P6").Value = 1.306431210178389
$ws.Range("Q6").Value = 0.5384892413853197
$ws.Range("R6").Value = 0.1824632656666271
$ws.Range("S6").Value = 1.275350349221521
$ws.Range("T6").Value = 0.239329294593088
$ws.Range("U6").Value = 1.282651499458469
$ws.Range("P7").Value = 0.01782614581474324
$ws.Range("Q7").Value = 1.001696046634265
$ws.Range("R7").Value = 0.2496393557079514
$ws.Range("S7").Value = 0.5608314417700987
$ws.Range("T7").Value = 0.004458362267418854
$ws.Range("U7").Value = 0.6951883024290437
$ws.Range("P8").Value = 1.307133394122553
$ws.Range("Q8").Value = 0.3985664986945726
$ws.Range("R8").Value = 0.1802111618208298
$ws.Range("S8").Value = 1.368412169729564
$ws.Range("T8").Value = 0.2374190859997568
$ws.Range("U8").Value = 1.300845597702918
$ws.Range("P9").Value = 0.1436693283071879
$ws.Range("Q9").Value = 0.3656952702091088
$ws.Range("R9").Value = 0.01857158169352735
$ws.Range("S9").Value = 0.2674942102791997
$ws.Range("T9").Value = 0.002664780436688747
$ws.Range("U9").Value = 0.3003132600935411
$ws.Range("P10").Value = 1.306120948283209
$ws.Range("Q10").Value = 0.5962974751618334
$ws.Range("R10").Value = 0.1790284476818534
$ws.Range("S10").Value = 1.346510245473203
$ws.Range("T10").Value = 0.2340797711025001
$ws.Range("U10").Value = 1.249551901643145
$ws.Range("P11").Value = 0.14328651417766
$ws.Range("Q11").Value = 0.3362269797239298
$ws.Range("R11").Value = 0.0009646390862442004
$ws.Range("S11").Value = 0.2548289831718226
$ws.Range("T11").Value = 0.0001387130000064089
$ws.Range("U11").Value = 0.2564755472855793
$ws.Range("P12").Value = 1.311097789796539
$ws.Range("Q12").Value = 0.4514431168005566
$ws.Range("R12").Value = 0.1793556356172296
$ws.Range("S12").Value = 1.296064864780603
$ws.Range("T12").Value = 0.2343681912243217
$ws.Range("U12").Value = 1.283193381054588
